$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.496.19'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.982.75'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'382.08"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'103.99"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('D7').Value = "'0.545"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('D10').Value = "'36.69"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('D11').Value = "'0.137"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').Value = "'0.0857"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '3.452.85'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = "'18.45"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = "'7.80"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').Value = '3.003.41'
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = "'0.998"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').Value = '51.519.58'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').Value = "'12.59"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').Value = '0.0₃0963'
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('D23').Value = "'70.37"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.30%  '
$ws.Range('D24').Value = "'267.03"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('D26').Value = "'7.85"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.78%  '
$ws.Range('D27').Value = "'7.33"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('D28').Value = "'0.168"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.48%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = "'26.11"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').Value = "'10.36"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.58%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = "'34.80"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.05%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = "'51.51"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').Value = "'2.07"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('D36').Value = "'0.0441"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').Value = "'3.28"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.92%  '
$ws.Range('D39').Value = "'16.78"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.33%  '
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('D41').Value = "'2.56"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('D42').Value = "'1.84"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.85%  '
$ws.Range('D43').Value = "'126.72"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.16%  '
$ws.Range('D44').Value = "'3.80"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.26%  '
$ws.Range('D45').Value = "'21.46"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').Value = "'0.271"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = "'2.37"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.49%  '
$ws.Range('D49').Value = '2.027.88'
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('D50').Value = "'0.0335"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('E51').Value = '  +15.40%  '
